# Fix errors in Excel exercises:
# The "Social security" row (12) on the Solution sheet was computed from the
# "Basic salary" row (6) instead of the "GROSS SALARY" row (10). Correct the
# formula so that Social security is calculated on the gross salary, and
# propagate it across the row the same way a user would (fix B12, then copy
# the corrected formula across C12:F12).

$wb = $excel.ActiveWorkbook
$wsSolution = $wb.Worksheets.Item("Solution")
$wsExercise = $wb.Worksheets.Item("Exercise")

$wsSolution.Range("B12").Formula = '=+B10*$B25'
$wsSolution.Range("C12:F12").Formula = '=+C10*$B25'

# Restore the view: scroll the Solution sheet back to the top (it had been
# left scrolled down to row 10) without leaving it as the active tab.
$wsSolution.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsExercise.Activate()
